$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Capture existing hyperlinks (address + target URL) before touching layout.
$links = @()
foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    $target = $h.Address()
    $links += , @($addr, $target)
}

# 2) Remove all existing hyperlinks (keeps cell text/style intact).
while ($ws.Hyperlinks.Count -gt 0) {
    foreach ($h in $ws.Hyperlinks) {
        $h.Delete()
        break
    }
}

# 3) Insert a new column before column B; this shifts old B -> C, old C -> D
#    (values & styles move with it).
$ws.Columns("B").Insert()

# 4) The newly inserted column B cells inherit the old column B formatting
#    (hyperlink style); clear that so they are plain/default-styled cells.
$ws.Range("B1:B6").ClearFormats()

# 5) Populate the new "Badge" column.
$ws.Range("B1").Value = "Badge"
$ws.Range("B2").Value = "Artifact Evaluated"
$ws.Range("B3").Value = "Artifact Evaluated"
$ws.Range("B4").Value = "Artifact Evaluated"
$ws.Range("B5").Value = "Artifact Evaluated"
$ws.Range("B6").Value = "Artifact Evaluated"

# 5b) Widen column A slightly and give the new column B an explicit width,
#     matching the saved workbook's column layout as closely as possible.
$ws.Columns("A").ColumnWidth = 59.666666666666664
$ws.Columns("B").ColumnWidth = 55.666666666666664

# 6) Re-create the hyperlinks, shifting any that were in column B or C one
#    column to the right (B->C, C->D); column A references stay the same.
foreach ($l in $links) {
    $addr = $l[0]
    $target = $l[1]
    $col = $addr.Substring(1, $addr.IndexOf('$', 1) - 1)
    $row = $addr.Substring($addr.IndexOf('$', 1) + 1)

    if ($col -eq "A") {
        $newAddr = "A" + $row
    } elseif ($col -eq "B") {
        $newAddr = "C" + $row
    } elseif ($col -eq "C") {
        $newAddr = "D" + $row
    }

    $ws.Hyperlinks.Add($ws.Range($newAddr), $target)
    # Adding a hyperlink via COM re-applies formatting as a direct override;
    # reapply the standard "Hyperlink" cell style so the cell matches the
    # original look (underlined themed link text) used throughout the sheet.
    $ws.Range($newAddr).Style = "Hyperlink"
}

# 7) Restore the selection seen in the saved workbook.
$ws.Range("B12").Select()
